$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 33
$prev = 32

# Values
$ws.Cells.Item($row, 1).Value = 32
$ws.Cells.Item($row, 2).Value = "gibraltar"
$ws.Cells.Item($row, 3).Value = "national-league"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45261.875
$ws.Cells.Item($row, 6).Value = "Lincoln Red Imps"
$ws.Cells.Item($row, 7).Value = 2
$ws.Cells.Item($row, 8).Value = "Lions Gibraltar"
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 1.02
$ws.Cells.Item($row, 11).Value = "01/12/2023 11:23"
$ws.Cells.Item($row, 12).Value = 1.02
$ws.Cells.Item($row, 13).Value = "01/12/2023 11:23"
$ws.Cells.Item($row, 14).Value = 18.12
$ws.Cells.Item($row, 15).Value = "01/12/2023 13:28"
$ws.Cells.Item($row, 16).Value = 18.12
$ws.Cells.Item($row, 17).Value = "01/12/2023 13:28"
$ws.Cells.Item($row, 18).Value = 22.84
$ws.Cells.Item($row, 19).Value = "01/12/2023 13:28"
$ws.Cells.Item($row, 20).Value = 22.84
$ws.Cells.Item($row, 21).Value = "01/12/2023 13:28"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/gibraltar/national-league/lincoln-red-imps-lions-gibraltar/C4qcQkaK/"

# Match existing formatting: column A (bordered/bold/centered index) and
# column E (date-time number format) reuse the same style as the row above.
$ws.Cells.Item($prev, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)

$ws.Cells.Item($prev, 5).Copy()
$ws.Cells.Item($row, 5).PasteSpecial(-4122)

$excel.CutCopyMode = $false
